$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '87.123.76'
$ws.Range('E2').Value = '  -2.93%  '
$ws.Range('D3').Value = '3.143.10'
$ws.Range('E3').Value = '  -6.91%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '203.82'
$ws.Range('E5').Value = '  -7.44%  '
$ws.Range('D6').Value = '603.70'
$ws.Range('E6').Value = '  -7.22%  '
$ws.Range('D7').Value = '0.372'
$ws.Range('E7').Value = '  -9.43%  '
$ws.Range('D8').Value = '0.658'
$ws.Range('E8').Value = '  +7.17%  '
$ws.Range('E9').Value = '  -0.04%  '
$ws.Range('D10').Value = '3.144.12'
$ws.Range('E10').Value = '  -6.89%  '
$ws.Range('D11').Value = '0.530'
$ws.Range('E11').Value = '  -12.12%  '
$ws.Range('E12').Value = '  +4.79%  '
$ws.Range('D13').Value = '0.0000241'
$ws.Range('E13').Value = '  -17.09%  '
$ws.Range('D14').Value = '3.725.00'
$ws.Range('E14').Value = '  -6.88%  '
$ws.Range('D15').Value = '5.21'
$ws.Range('E15').Value = '  -6.90%  '
$ws.Range('D16').Value = '86.754.19'
$ws.Range('E16').Value = '  -3.22%  '
$ws.Range('D17').Value = '31.84'
$ws.Range('E17').Value = '  -13.94%  '
$ws.Range('D18').Value = '3.181.47'
$ws.Range('E18').Value = '  -5.56%  '
$ws.Range('D19').Value = '2.98'
$ws.Range('E19').Value = '  -6.56%  '
$ws.Range('D20').Value = '13.30'
$ws.Range('E20').Value = '  -10.25%  '
$ws.Range('D21').Value = '411.19'
$ws.Range('E21').Value = '  -10.61%  '
$ws.Range('D22').Value = '8.43'
$ws.Range('E22').Value = '  -12.90%  '
$ws.Range('D23').Value = '5.07'
$ws.Range('E23').Value = '  -8.48%  '
$ws.Range('D24').Value = '5.12'
$ws.Range('E24').Value = '  -8.63%  '
$ws.Range('D25').Value = '11.81'
$ws.Range('E25').Value = '  -8.46%  '
$ws.Range('D26').Value = '3.309.14'
$ws.Range('E26').Value = '  -5.42%  '
$ws.Range('D27').Value = '73.10'
$ws.Range('E27').Value = '  -7.34%  '
$ws.Range('E28').Value = '  -10.75%  '
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +0.19%  '
$ws.Range('D30').Value = '0.163'
$ws.Range('E30').Value = '  -20.12%  '
$ws.Range('E31').Value = '  +0.19%  '
$ws.Range('D32').Value = '534.34'
$ws.Range('E32').Value = '  -10.76%  '
$ws.Range('D33').Value = '8.22'
$ws.Range('E33').Value = '  -12.34%  '
$ws.Range('D34').Value = '1.29'
$ws.Range('E34').Value = '  -19.13%  '
$ws.Range('D35').Value = '1.83'
$ws.Range('E35').Value = '  -13.34%  '
$ws.Range('D36').Value = '6.56'
$ws.Range('E36').Value = '  -12.03%  '
$ws.Range('D37').Value = '0.132'
$ws.Range('E37').Value = '  -8.42%  '
$ws.Range('D38').Value = '21.76'
$ws.Range('E38').Value = '  -7.35%  '
$ws.Range('D39').Value = '1.00'
$ws.Range('E39').Value = '  +0.12%  '
$ws.Range('D40').Value = '21.78'
$ws.Range('E40').Value = '  -0.54%  '
$ws.Range('E41').Value = '  -7.96%  '
$ws.Range('E42').Value = '  -0.01%  '
$ws.Range('E43').Value = '  -13.45%  '
$ws.Range('D44').Value = '0.369'
$ws.Range('E44').Value = '  -13.85%  '
$ws.Range('D45').Value = '149.07'
$ws.Range('E45').Value = '  -5.65%  '
$ws.Range('D46').Value = '171.05'
$ws.Range('E46').Value = '  -9.97%  '
$ws.Range('D47').Value = '42.94'
$ws.Range('E47').Value = '  -7.39%  '
$ws.Range('E48').Value = '  +5.60%  '
$ws.Range('D49').Value = '1.25'
$ws.Range('E49').Value = '  -15.98%  '
$ws.Range('D50').Value = '3.95'
$ws.Range('E50').Value = '  -12.55%  '
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').Value = '0.690'
$ws.Range('E51').Value = '  -12.70%  '
